$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# "Periodo Mora" labels shift down one row (E16:E22), and the "Valor Mora" (F)
# and "Salario Basico" (G) figures for the new period table are updated.
$ws.Range("E16").Value = "1906"
$ws.Range("E17").Value = "1907"
$ws.Range("E18").Value = "1908"
$ws.Range("E19").Value = "1909"
$ws.Range("E20").Value = "1910"
$ws.Range("E21").Value = "1912"
$ws.Range("E22").Value = "2002"

$ws.Range("F16").Value = 32000
$ws.Range("F17").Value = 32000
$ws.Range("F18").Value = 32000
$ws.Range("F19").Value = 32000
$ws.Range("F20").Value = 33125
$ws.Range("F21").Value = 33125
$ws.Range("F22").Value = 18134

$ws.Range("G16").Value = 800000
$ws.Range("G17").Value = 800000
$ws.Range("G18").Value = 800000
$ws.Range("G19").Value = 800000
$ws.Range("G20").Value = 800000
$ws.Range("G21").Value = 800000
$ws.Range("G22").Value = 800000
